$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("REPORT")
$ws.Activate()

# Capture column A's display width before inserting, so the new column B
# (which Excel insert naturally copies the left-neighbour's format into)
# ends up the same visual width as column A.
$colAWidth = $ws.Columns("A").ColumnWidth

# Insert a new blank column before column B - this shifts the existing
# B (coefficient values) and C columns one to the right (-> C and D).
$ws.Columns("B").Insert() | Out-Null

# Match column B's width to column A's.
$ws.Columns("B").ColumnWidth = $colAWidth

# Put a literal "=" (quote-prefixed text, not a formula) in B2:B7 next to
# each of the six affine coefficients.
$ws.Range("B2").Value = "'="
$ws.Range("B3").Value = "'="
$ws.Range("B4").Value = "'="
$ws.Range("B5").Value = "'="
$ws.Range("B6").Value = "'="
$ws.Range("B7").Value = "'="

# Update the visible selection to match the example being highlighted.
$ws.Range("A2:C7").Select() | Out-Null
